# Apply edits matching commit: "completed the tests of client.py"
# - fill in the Client class unit-test rows (Method Inputs / Expected Result)
# - widen/narrow a few columns and re-size several rows to fit the new text
# - update the selection to the cell the author finished on

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in new Client test-case data (columns F and G, rows 7-16) ---
$ws.Range("F7").Value = 'client_number=123,             first_name="Lily",      last_name="Green",        email_address="lilygreen@gmail.com"'
$ws.Range("G7").Value = 'The Client instance should be created successfully with the attributes correctly set.'
$ws.Range("F8").Value = 'client_number=None  first_name="Lily",      last_name="Green",        email_address="lilygreen@gmail.com"'
$ws.Range("G8").Value = 'Raises ValueError successfully'
$ws.Range("F9").Value = 'client_number=123,             first_name="  ",      last_name="Green",        email_address="lilygreen@gmail.com"'
$ws.Range("G9").Value = 'Raises ValueError successfully'
$ws.Range("F10").Value = 'client_number=123,             first_name="Lily",      last_name="  ",        email_address="lilygreen@gmail.com"'
$ws.Range("G10").Value = 'Raises ValueError successfully'
$ws.Range("F11").Value = 'client_number=123,             first_name="Lily",      last_name="Green",        email_address="lilygreengmail"'
$ws.Range("G11").Value = 'reset the email address attribute to "email@pixell-river.com" succesfully'
$ws.Range("F12").Value = 'client_number=123'
$ws.Range("G12").Value = 'succesfully return client_number attribute'
$ws.Range("F13").Value = 'first_name="Lily"'
$ws.Range("G13").Value = 'succesfully return  first_name attribute'
$ws.Range("F14").Value = 'last_name="Green"'
$ws.Range("G14").Value = 'succesfully return last_name attribute'
$ws.Range("F15").Value = 'email_address="lilygreengmail"'
$ws.Range("G15").Value = 'succesfully return email_address attribute'
$ws.Range("F16").Value = 'self.client = Client(123,"Lily","Green","lilygreen@gmail.com")'
$ws.Range("G16").Value = 'succesfully return a string as setted'

# G16 picks up the bold "answered" style used by the rest of the table
# (F16 intentionally keeps its original, non-bold style)
$ws.Range("G16").Font.Bold = $true

# --- Row heights: rows 7-12 grow to fit the newly entered, wrapped text ---
$ws.Rows.Item(7).RowHeight = 75.75
$ws.Rows.Item(8).RowHeight = 78.75
$ws.Rows.Item(9).RowHeight = 73.5
$ws.Rows.Item(10).RowHeight = 70.5
$ws.Rows.Item(11).RowHeight = 72
$ws.Rows.Item(12).RowHeight = 39.4

# --- Column widths: minor manual resize of the data columns ---
$ws.Columns.Item("B").ColumnWidth = 11.833333333333334
$ws.Columns.Item("C").ColumnWidth = 9.333333333333334
$ws.Columns.Item("D").ColumnWidth = 23.166666666666668
$ws.Columns.Item("E").ColumnWidth = 9.5
$ws.Columns.Item("F").ColumnWidth = 22.333333333333332
$ws.Columns.Item("G").ColumnWidth = 25.5

# --- Selection left on K14 after the edits ---
$ws.Range("K14").Select()

